$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# --- Row 6: weight of finished product (plain value, no formula) ---
$ws.Range("I6").Value = 468

# --- Row 7: weight = SUM(G7*H7) ---
$ws.Range("I7").Formula = "=SUM(G7*H7)"

# --- Rows 8,10,11,12,13: fill in per-part weights (column G) ---
$ws.Range("G8").Value = 10
$ws.Range("G10").Value = 12
$ws.Range("G11").Value = 4
$ws.Range("G12").Value = 4
$ws.Range("G13").Value = 4

# --- Rows 8-16: shared weight formula SUM(G{row}*H{row}) in column I ---
$ws.Range("I8:I16").Formula = "=SUM(G8*H8)"

# --- Row 15: extra helper formula ---
$ws.Range("L15").Formula = "=56*2"

# --- Row 17: plain weight value (no formula) ---
$ws.Range("I17").Value = 12

# --- Rows 19-20: replace computed weight formulas with fixed values ---
$ws.Range("I19").Value = 91
$ws.Range("I20").Value = 80

# --- Row 21: clear the old SUM(I19:I20) formula ---
$ws.Range("I21").ClearContents()

# --- Row 26: drop the old G26 total, add new I26 total + J26 note ---
$ws.Range("G26").ClearContents()
$ws.Range("I26").Formula = "=SUM(I4:I25)"
$ws.Range("J26").Value = "무게 1kg 이내로"

# --- Update the view: scroll position + selected cell ---
$ws.Range("A7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select() | Out-Null
